# Updated cryptos list on Sun Oct 29 21:16:20 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep a numeric-looking string as TEXT (matches how
    # the source data feed writes these "Price" values as plain strings,
    # not floats) instead of letting Excel auto-coerce it to a number.
    $ws.Range($range).NumberFormat = "@"
    $ws.Range($range).Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "34.615.10"
$ws.Range("E2").Value = "  +1.13%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.802.39"
$ws.Range("E3").Value = "  +1.14%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.14%  "

# Row 5 - BNB
Set-TextValue "D5" "227.48"
$ws.Range("E5").Value = "  +0.49%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +2.01%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.16%  "

# Row 8 - Solana
Set-TextValue "D8" "32.85"
$ws.Range("E8").Value = "  +3.00%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +1.62%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.0699"
$ws.Range("E10").Value = "  +1.19%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0951"
$ws.Range("E11").Value = "  +0.39%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "2.059.30"
$ws.Range("E12").Value = "  +0.97%  "

# Rows 13 and 14 swap order: Chainlink <-> WrappedEther
# Row 13 becomes WrappedEther
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.810.49"
$ws.Range("E13").Value = "  +1.85%  "

# Row 14 becomes Chainlink
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D14" "11.14"
$ws.Range("E14").Value = "  +1.09%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +2.70%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "34.604.71"
$ws.Range("E16").Value = "  +1.31%  "

# Row 17 - Polkadot
$ws.Range("E17").Value = "  +3.04%  "

# Row 18 - Litecoin
Set-TextValue "D18" "69.01"
$ws.Range("E18").Value = "  +1.64%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.0₃0806"
$ws.Range("E19").Value = "  +0.91%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "247.76"
$ws.Range("E20").Value = "  +0.25%  "

# Row 21 - Avalanche
Set-TextValue "D21" "11.31"
$ws.Range("E21").Value = "  +2.96%  "

# Row 22 - Dai (unchanged)

# Row 23 - Uniswap
Set-TextValue "D23" "4.19"
$ws.Range("E23").Value = "  +2.24%  "

# Row 24 - Monero
Set-TextValue "D24" "167.26"
$ws.Range("E24").Value = "  +2.99%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  +1.27%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  +1.60%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "16.63"
$ws.Range("E27").Value = "  +1.97%  "

# Row 28 - Stellar
Set-TextValue "D28" "0.117"
$ws.Range("E28").Value = "  +2.34%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  -0.15%  "

# Row 30 - InternetComputer(DFINITY)
Set-TextValue "D30" "4.11"
$ws.Range("E30").Value = "  +11.02%  "

# Rows 31 and 32 swap order: PancakeSwap <-> Hedera
# Row 31 becomes Hedera
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D31" "0.0526"
$ws.Range("E31").Value = "  +1.11%  "

# Row 32 becomes PancakeSwap
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D32" "1.24"
$ws.Range("E32").Value = "  +0.57%  "

# Row 33 - Filecoin
Set-TextValue "D33" "3.82"
$ws.Range("E33").Value = "  +1.97%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +2.80%  "

# Row 35 - Maker
Set-TextValue "D35" "1.432.87"
$ws.Range("E35").Value = "  -0.81%  "

# Row 36 - RenderToken
Set-TextValue "D36" "2.59"
$ws.Range("E36").Value = "  +7.27%  "

# Row 37 - ImmutableX
Set-TextValue "D37" "0.673"
$ws.Range("E37").Value = "  +3.12%  "

# Row 38 - TrustWalletToken
$ws.Range("E38").Value = "  +2.50%  "

# Row 39 - VeChain
Set-TextValue "D39" "0.0193"
$ws.Range("E39").Value = "  +0.56%  "

# Row 40 - Aave
Set-TextValue "D40" "85.60"
$ws.Range("E40").Value = "  +6.71%  "

# Row 41 - HuobiToken
$ws.Range("E41").Value = "  +1.09%  "

# Row 42 - ARBITRUM
Set-TextValue "D42" "0.943"
$ws.Range("E42").Value = "  +2.11%  "

# Row 43 - MXToken
$ws.Range("E43").Value = "  +3.39%  "

# Row 44 - InjectiveProtocol
Set-TextValue "D44" "13.75"
$ws.Range("E44").Value = "  +0.83%  "

# Row 45 - Kaspa
Set-TextValue "D45" "0.0527"
$ws.Range("E45").Value = "  +3.49%  "

# Row 46 - FraxShare
Set-TextValue "D46" "6.11"
$ws.Range("E46").Value = "  +0.62%  "

# Row 47 - WEMIXToken
$ws.Range("E47").Value = "  +0.25%  "

# Row 48 - RocketPoolETH
Set-TextValue "D48" "1.958.41"
$ws.Range("E48").Value = "  +0.87%  "

# Row 49 - Quant
Set-TextValue "D49" "106.22"
$ws.Range("E49").Value = "  +1.38%  "

# Row 50 - PaxDollar
$ws.Range("E50").Value = "  -0.19%  "

# Row 51 - BabyDogeCoin
$ws.Range("E51").Value = "  -6.17%  "
